$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-10 (A:T), reflecting the expanded
# ECs/FAPs/sCs x ECs/FAPs/sCs Angpt1->Tek NATMI LR-pair cross table.
$data = @(
  ,@("ECs","Angpt1","Tek","ECs",1,0.3333333333333333,0.09342833333333334,0.280285,0.004971127078137211,0.004971127078137211,3,1,39.413899,118.241697,0.7310109930444597,0.7310109930444597,3.682374893738333,33.141374043645,0.003633948541939286,0.003633948541939286)
  ,@("ECs","Angpt1","Tek","FAPs",1,0.3333333333333333,0.09342833333333334,0.280285,0.004971127078137211,0.004971127078137211,3,1,11.40791933333333,34.223758,0.2115830874897988,0.2115830874897987,1.065822890114444,9.59240601103,0.001051806415496413,0.001051806415496413)
  ,@("ECs","Angpt1","Tek","sCs",1,0.3333333333333333,0.09342833333333334,0.280285,0.004971127078137211,0.004971127078137211,3,1,3.095153333333334,9.28546,0.05740591946574151,0.0574059194657415,0.2891750173444445,2.6025751561,0.0002853721207015117,0.0002853721207015116)
  ,@("FAPs","Angpt1","Tek","ECs",3,1,14.289121,42.867363,0.760294375288143,0.760294375288143,3,1,39.413899,118.241697,0.7310109930444597,0.7310109930444597,563.189971892779,5068.709747035011,0.5557835462855025,0.5557835462855025)
  ,@("FAPs","Angpt1","Tek","FAPs",3,1,14.289121,42.867363,0.760294375288143,0.760294375288143,3,1,11.40791933333333,34.223758,0.2115830874897988,0.2115830874897987,163.0091397122393,1467.082257410154,0.1608654313245931,0.1608654313245931)
  ,@("FAPs","Angpt1","Tek","sCs",3,1,14.289121,42.867363,0.760294375288143,0.760294375288143,3,1,3.095153333333334,9.28546,0.05740591946574151,0.0574059194657415,44.22702049355334,398.04318444198,0.04364539767804739,0.04364539767804738)
  ,@("sCs","Angpt1","Tek","ECs",3,1,4.411646,13.234938,0.2347344976337197,0.2347344976337198,3,1,39.413899,118.241697,0.7310109930444597,0.7310109930444597,173.880169867754,1564.921528809786,0.1715934982170178,0.1715934982170179)
  ,@("sCs","Angpt1","Tek","FAPs",3,1,4.411646,13.234938,0.2347344976337197,0.2347344976337198,3,1,11.40791933333333,34.223758,0.2115830874897988,0.2115830874897987,50.32770169522267,452.9493152570041,0.04966584974970928,0.04966584974970929)
  ,@("sCs","Angpt1","Tek","sCs",3,1,4.411646,13.234938,0.2347344976337197,0.2347344976337198,3,1,3.095153333333334,9.28546,0.05740591946574151,0.0574059194657415,13.65472082238667,122.89248740148,0.01347514966699261,0.01347514966699261)
)

$r = 2
foreach ($row in $data) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}
